$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.406.68"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.39"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.56"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7045"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07922"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3133"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07813"
$ws.Range("E11").Value = "  -4.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.893.86"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.68"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.172"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7030"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.527"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008429"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.450.57"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.146.89"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.674"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.014"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.62"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.506"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.313"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.259"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.215"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05265"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.182"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7534"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01877"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.279.80"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.770"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8975"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.19"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.018"
$ws.Range("E43").Value = "  -6.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.94"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.044.62"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000127"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.607"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5198"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4298"
$ws.Range("E51").Value = "  -1.28%  "
